$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank rows before the "Interpretations" block (old row 18),
#    pushing Interpretations (18-21) -> (20-23) and Runs (26-30) -> (28-32).
$ws.Rows("18:19").Insert()

# 2. New "Users:" table
$ws.Range("A35").Value = "Users:"
$ws.Range("B35").Value = "username"
$ws.Range("B36").Value = "password"
$ws.Range("C36").Value = "varchar"
$ws.Range("C35").Value = "varchar"

# 3. New "Panels:" table (and the matching Fkey reference back on the
#    Patient_info/Panel row)
$ws.Range("A39").Value = "Panels:"
$ws.Range("C7").Value = "Fkey - panels"
$ws.Range("D7").Value = "X"
$ws.Range("B39").Value = "Type panel"
$ws.Range("C39").Value = "varchar"
$ws.Range("D39").Value = "pv2-1, exom, filtex"

# 4. Rename a few fields
$ws.Range("B2").Value = "Patient_ID"
$ws.Range("B4").Value = "Clinical_info"
$ws.Range("A9").Value = "raw_variants"

# 5. New genome-build row in raw_variants, and PK marker on Runs/SBS
$ws.Range("B14").Value = "hg"
$ws.Range("C14").Value = "varchar"
$ws.Range("D14").Value = "hg19"
$ws.Range("D28").Value = "PK"

# 6. A couple of dtype fill-ins that reuse already-existing strings
$ws.Range("C6").Value = "varchar"
$ws.Range("C12").Value = "Varchar"
$ws.Range("C13").Value = "Varchar"

# 7. The old QC fields (Mean_target_cov / % X > 20 X / % X > 30 X), now shifted
#    to rows 30-32 by the row insert above, move down into their own "QC:"
#    block at rows 42-45.
$ws.Range("A30:F32").ClearContents()

$ws.Range("A42").Value = "QC:"
$ws.Range("B42").Value = "Mean_target_cov"
$ws.Range("C42").Value = "Float"
$ws.Range("B43").Value = "% X > 20 X"
$ws.Range("C43").Value = "Float"
$ws.Range("B44").Value = "% X > 30 X"
$ws.Range("C44").Value = "Float"
$ws.Range("B45").Value = "info fra Hsmetrics (JSON?)"

# 8. Bold the new section headers to match the "Table name" header style
$ws.Range("A39").Font.Bold = $true
$ws.Range("A42").Font.Bold = $true

# 9. Column C is now used for longer dtype text (e.g. "Fkey - panels"); widen it
$ws.Columns("C").ColumnWidth = 16.1640625

# 10. Restore the selection roughly where the author left off editing
$ws.Range("F32").Select()

# 11. Page setup matches the exported workbook (A4, portrait print)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
